# Generate Report for Archive
#
# 1. Status text "Ready for handoff" -> "In Translation" (appears in the
#    "Status" column on every sheet: Overview!E2:F4, zh-cn!C2:C4, de-de!C2:C4).
# 2. Narrow the "Status" column(s) from ~17.2 to ~13.4 characters wide:
#    Overview columns E & F, and column C on the zh-cn / de-de sheets.

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    foreach ($cell in $used.Cells) {
        # NOTE: keep the string literal on the LEFT of -eq. PowerShell's -eq
        # coerces the right-hand side to the left-hand side's type, so with a
        # boolean cell value (e.g. "True"/"False" status flags) on the left,
        # "$true -eq 'Ready for handoff'" would wrongly evaluate to $true.
        if ("Ready for handoff" -eq $cell.Value()) {
            $cell.Value = "In Translation"
        }
    }
}

# Narrower "Status" columns.
# (The underlying engine snaps COM ColumnWidth to an internal 1/6-character
# grid, so 12.5 is the nearest settable value to the recorded raw width of
# 13.4101848602295 characters.)
$narrowWidth = 12.5

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E1").ColumnWidth = $narrowWidth
$wsOverview.Range("F1").ColumnWidth = $narrowWidth

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C1").ColumnWidth = $narrowWidth

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C1").ColumnWidth = $narrowWidth
